{"js": "// Update the worksheet date title and every arithmetic-problem cell in the\n// single table, matching the target revision (new date + 100 new problems).\n// Each table cell keeps its existing run formatting because we only change\n// the cell text via Table.values (same technique Office.js uses to rewrite\n// cell contents in place) and the title run via Paragraph.insertText with\n// InsertLocation.replace (keeps the run's rFonts/sz formatting intact).\n\nconst NEW_TITLE = \"2025-02-17 Monday\";\n\n// New value for every cell, in row-major order (20 rows x 5 cols = 100).\nconst NEW_ROWS = [\n  [\"24+24=\", \"3+40=\", \"56+11=\", \"38+53=\", \"27+33=\"],\n  [\"96-82=\", \"32+30=\", \"44-16=\", \"45+25=\", \"60-4=\"],\n  [\"19-3=\", \"74-58=\", \"6+70=\", \"34-20=\", \"55-10=\"],\n  [\"82-68=\", \"80-57=\", \"17-13=\", \"18+6=\", \"97-17=\"],\n  [\"28-18=\", \"1+50=\", \"88-80=\", \"13+31=\", \"80-29=\"],\n  [\"77-19=\", \"89-87=\", \"1+61=\", \"27+43=\", \"70-42=\"],\n  [\"30-18=\", \"58+37=\", \"8+46=\", \"12+32=\", \"99-86=\"],\n  [\"83-35=\", \"6+53=\", \"75-66=\", \"27+10=\", \"52+38=\"],\n  [\"21+4=\", \"85-69=\", \"61+14=\", \"33-2=\", \"31+53=\"],\n  [\"43+6=\", \"55-0=\", \"8+64=\", \"11+36=\", \"42-36=\"],\n  [\"94-2=\", \"89-54=\", \"68-20=\", \"62-61=\", \"58-23=\"],\n  [\"3+23=\", \"67-58=\", \"87+6=\", \"27+67=\", \"83-71=\"],\n  [\"68-40=\", \"69-58=\", \"99-60=\", \"16+82=\", \"94-62=\"],\n  [\"86-14=\", \"67-58=\", \"55+11=\", \"39-26=\", \"54+21=\"],\n  [\"10+86=\", \"95-74=\", \"61+14=\", \"68-53=\", \"73+22=\"],\n  [\"97-49=\", \"68-56=\", \"9-8=\", \"61-55=\", \"10+3=\"],\n  [\"11+75=\", \"37-8=\", \"7+38=\", \"81-9=\", \"41+10=\"],\n  [\"61+20=\", \"71-17=\", \"19+69=\", \"71-24=\", \"77-54=\"],\n  [\"4+1=\", \"32+1=\", \"18+55=\", \"45+5=\", \"21+8=\"],\n  [\"70+22=\", \"89-20=\", \"74-5=\", \"63-61=\", \"13+16=\"],\n];\n\n// 1) Update the date/weekday title paragraph (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(NEW_TITLE, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Update every cell of the (single) table in one shot using `.values`,\n// which rewrites cell text while leaving cell/run formatting untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = NEW_ROWS;\nawait context.sync();\n", "ps1": "# Update the worksheet date title and every arithmetic-problem cell in the\n# single table, matching the target revision (new date + 100 new problems).\n# Per-cell Range.Text assignment swaps only the text run content, so the\n# existing cell/run formatting (TimeNewRoman font, sz 30, left-justified\n# paragraph, etc.) is left untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday title (first paragraph of the document).\n$d.Paragraphs.Item(1).Range.Text = \"2025-02-17 Monday\"\n\n# 2) Update every cell of the (single) table, row by row / column by column.\n$newValues = @(\n    @(\"24+24=\", \"3+40=\", \"56+11=\", \"38+53=\", \"27+33=\"),\n    @(\"96-82=\", \"32+30=\", \"44-16=\", \"45+25=\", \"60-4=\"),\n    @(\"19-3=\", \"74-58=\", \"6+70=\", \"34-20=\", \"55-10=\"),\n    @(\"82-68=\", \"80-57=\", \"17-13=\", \"18+6=\", \"97-17=\"),\n    @(\"28-18=\", \"1+50=\", \"88-80=\", \"13+31=\", \"80-29=\"),\n    @(\"77-19=\", \"89-87=\", \"1+61=\", \"27+43=\", \"70-42=\"),\n    @(\"30-18=\", \"58+37=\", \"8+46=\", \"12+32=\", \"99-86=\"),\n    @(\"83-35=\", \"6+53=\", \"75-66=\", \"27+10=\", \"52+38=\"),\n    @(\"21+4=\", \"85-69=\", \"61+14=\", \"33-2=\", \"31+53=\"),\n    @(\"43+6=\", \"55-0=\", \"8+64=\", \"11+36=\", \"42-36=\"),\n    @(\"94-2=\", \"89-54=\", \"68-20=\", \"62-61=\", \"58-23=\"),\n    @(\"3+23=\", \"67-58=\", \"87+6=\", \"27+67=\", \"83-71=\"),\n    @(\"68-40=\", \"69-58=\", \"99-60=\", \"16+82=\", \"94-62=\"),\n    @(\"86-14=\", \"67-58=\", \"55+11=\", \"39-26=\", \"54+21=\"),\n    @(\"10+86=\", \"95-74=\", \"61+14=\", \"68-53=\", \"73+22=\"),\n    @(\"97-49=\", \"68-56=\", \"9-8=\", \"61-55=\", \"10+3=\"),\n    @(\"11+75=\", \"37-8=\", \"7+38=\", \"81-9=\", \"41+10=\"),\n    @(\"61+20=\", \"71-17=\", \"19+69=\", \"71-24=\", \"77-54=\"),\n    @(\"4+1=\", \"32+1=\", \"18+55=\", \"45+5=\", \"21+8=\"),\n    @(\"70+22=\", \"89-20=\", \"74-5=\", \"63-61=\", \"13+16=\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
